$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to remain text
# (otherwise Excel auto-converts "1.00" -> 1, "591.51" -> 591.51 as a number, etc.)
$textCells = @("D5", "D6", "D10", "D12", "D14", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D29", "D30", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell (row order matches the diff)
$ws.Range("D2").Value = "67.382.11"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "2.608.71"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "591.51"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "150.59"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "2.608.02"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  -3.22%  "
$ws.Range("D14").Value = "27.29"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "3.081.62"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").Value = "67.214.49"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "2.610.13"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "371.42"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "7.36"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").Value = "4.84"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  -3.81%  "
$ws.Range("D25").Value = "73.31"
$ws.Range("E25").Value = "  +4.63%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "2.744.27"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "578.25"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "0.0₃0987"
$ws.Range("E31").Value = "  -6.09%  "
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -4.08%  "
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").Value = "157.92"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").Value = "19.07"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("D40").Value = "1.86"
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("D44").Value = "17.13"
$ws.Range("E44").Value = "  +4.24%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "153.36"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("D47").Value = "0.0₆0284"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "1.68"
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0778"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").Value = "21.36"
$ws.Range("E51").Value = "  +1.23%  "
